# Update crypto price/volume data per the Sep 3 2023 GitHub Actions scrape run.
# Several rows also have their Coin/Link reordered/swapped (rank churn in the source feed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces each cell to stay plain text (prices like "1.008" or
# "6.070" would otherwise be auto-coerced to numbers and lose formatting/precision).
$ws.Range("D2").Value = "'25.900.66"
$ws.Range("E2").Value = "'  +0.27%  "
$ws.Range("D3").Value = "'1.647.89"
$ws.Range("E3").Value = "'  +0.84%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "'  +0.66%  "
$ws.Range("D5").Value = "'215.56"
$ws.Range("E5").Value = "'  -0.01%  "
$ws.Range("E6").Value = "'  +0.84%  "
$ws.Range("E7").Value = "'  +0.52%  "
$ws.Range("D8").Value = "'0.2574"
$ws.Range("E8").Value = "'  -0.09%  "
$ws.Range("D9").Value = "'0.06419"
$ws.Range("E9").Value = "'  -0.03%  "
$ws.Range("E10").Value = "'  -0.07%  "
$ws.Range("D11").Value = "'0.07783"
$ws.Range("E11").Value = "'  +1.14%  "
$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.691.95"
$ws.Range("E12").Value = "'  +3.49%  "
$ws.Range("B13").Value = "'Polkadot"
$ws.Range("C13").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.315"
$ws.Range("E13").Value = "'  +1.45%  "
$ws.Range("D14").Value = "'0.5469"
$ws.Range("E14").Value = "'  +0.16%  "
$ws.Range("D15").Value = "'0.0₅7903"
$ws.Range("E15").Value = "'  -0.44%  "
$ws.Range("D16").Value = "'65.08"
$ws.Range("E16").Value = "'  +2.41%  "
$ws.Range("D17").Value = "'26.000.01"
$ws.Range("E17").Value = "'  +0.63%  "
$ws.Range("D18").Value = "'1.007"
$ws.Range("E18").Value = "'  +0.56%  "
$ws.Range("D19").Value = "'196.92"
$ws.Range("E19").Value = "'  -3.15%  "
$ws.Range("D20").Value = "'4.429"
$ws.Range("E20").Value = "'  +2.42%  "
$ws.Range("D21").Value = "'10.04"
$ws.Range("E21").Value = "'  +0.92%  "
$ws.Range("D22").Value = "'6.070"
$ws.Range("E22").Value = "'  +1.65%  "
$ws.Range("D23").Value = "'1.009"
$ws.Range("E23").Value = "'  +0.62%  "
$ws.Range("D24").Value = "'1.859"
$ws.Range("E24").Value = "'  -3.13%  "
$ws.Range("D25").Value = "'141.11"
$ws.Range("E25").Value = "'  -0.05%  "
$ws.Range("D26").Value = "'0.1145"
$ws.Range("D27").Value = "'6.901"
$ws.Range("E27").Value = "'  +2.77%  "
$ws.Range("E28").Value = "'  +0.41%  "
$ws.Range("D29").Value = "'1.242"
$ws.Range("E29").Value = "'  -0.01%  "
$ws.Range("D30").Value = "'0.05028"
$ws.Range("E30").Value = "'  +0.00%  "
$ws.Range("D31").Value = "'3.273"
$ws.Range("E31").Value = "'  +0.08%  "
$ws.Range("D32").Value = "'3.200"
$ws.Range("E32").Value = "'  +0.41%  "
$ws.Range("D33").Value = "'1.544"
$ws.Range("E33").Value = "'  +0.50%  "
$ws.Range("D34").Value = "'2.371"
$ws.Range("E34").Value = "'  +0.52%  "
$ws.Range("D35").Value = "'0.8947"
$ws.Range("E35").Value = "'  -0.02%  "
$ws.Range("D36").Value = "'2.600"
$ws.Range("E36").Value = "'  -0.16%  "
$ws.Range("D37").Value = "'0.5552"
$ws.Range("E37").Value = "'  -1.02%  "
$ws.Range("D38").Value = "'1.132.51"
$ws.Range("E38").Value = "'  -3.59%  "
$ws.Range("D39").Value = "'0.01564"
$ws.Range("E39").Value = "'  +0.22%  "
$ws.Range("E40").Value = "'  +0.66%  "
$ws.Range("D41").Value = "'5.665"
$ws.Range("E41").Value = "'  -0.05%  "
$ws.Range("D42").Value = "'0.8159"
$ws.Range("E42").Value = "'  +1.15%  "
$ws.Range("D43").Value = "'99.75"
$ws.Range("E43").Value = "'  +0.27%  "
$ws.Range("D44").Value = "'0.0₈124"
$ws.Range("E44").Value = "'  +7.44%  "
$ws.Range("D45").Value = "'1.784.55"
$ws.Range("E45").Value = "'  +0.77%  "
$ws.Range("D46").Value = "'0.4542"
$ws.Range("E46").Value = "'  +0.69%  "
$ws.Range("B47").Value = "'Aave"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'55.31"
$ws.Range("E47").Value = "'  +0.76%  "
$ws.Range("B48").Value = "'Frax"
$ws.Range("C48").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.007"
$ws.Range("E48").Value = "'  +0.40%  "
$ws.Range("D49").Value = "'0.05097"
$ws.Range("E49").Value = "'  +1.09%  "
$ws.Range("B50").Value = "'USDD"
$ws.Range("C50").Value = "'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "'1.009"
$ws.Range("E50").Value = "'  +0.51%  "
$ws.Range("B51").Value = "'Algorand"
$ws.Range("C51").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.09554"
$ws.Range("E51").Value = "'  +2.88%  "
